$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1238.1207
$ws.Range("J17").Value = 1251.614
$ws.Range("L17").Value = 3754.842
$ws.Range("N17").Value = -4090.842

$ws.Range("H43").Value = 4606.7915
$ws.Range("I43").Value = 4148.8335
$ws.Range("J43").Value = 4759.4443
$ws.Range("K43").Value = 4148.8335
$ws.Range("L43").Value = 4759.4443
$ws.Range("M43").Value = -4079.8335
$ws.Range("N43").Value = -4897.4443

$ws.Range("H64").Value = 7416.6387
$ws.Range("I64").Value = 4400
$ws.Range("J64").Value = 7903.1934
$ws.Range("K64").Value = 4400
$ws.Range("L64").Value = 7903.1934
$ws.Range("M64").Value = -4152
$ws.Range("N64").Value = -8399.1934

$ws.Range("H67").Value = 7416.6387
$ws.Range("I67").Value = 4400
$ws.Range("J67").Value = 7903.1934
$ws.Range("K67").Value = 4400
$ws.Range("L67").Value = 7903.1934
$ws.Range("M67").Value = -3542
$ws.Range("N67").Value = -9619.1934

$ws.Range("H93").Value = 100000
$ws.Range("J93").Value = 100000
$ws.Range("L93").Value = 100000
$ws.Range("N93").Value = -104992

$ws.Range("H100").Value = 5776.353
$ws.Range("J100").Value = 7899.8
$ws.Range("L100").Value = 7899.8
$ws.Range("N100").Value = -8981.799999999999

$ws.Range("H116").Value = 5089.8
$ws.Range("I116").Value = 5089.8
$ws.Range("K116").Value = 5089.8
$ws.Range("M116").Value = -1647.8

$ws.Range("H137").Value = 1295
$ws.Range("I137").Value = 1283
$ws.Range("K137").Value = 3849
$ws.Range("M137").Value = -1299

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 105.35294
$ws.Range("I5").Value = 42.18182
$ws.Range("J5").Value = 221.16667
$ws.Range("K5").Value = 42.18182
$ws.Range("L5").Value = 221.16667
$ws.Range("M5").Value = 69.81818
$ws.Range("N5").Value = -445.16667

$ws.Range("H45").Value = 2104.8
$ws.Range("I45").Value = 1877.5
$ws.Range("K45").Value = 1877.5
$ws.Range("M45").Value = -1500.5

$ws.Range("H61").Value = 3461.44
$ws.Range("I61").Value = 3461.44
$ws.Range("K61").Value = 3461.44
$ws.Range("M61").Value = -3249.44

$ws.Range("H74").Value = 2153.2222
$ws.Range("I74").Value = 2190.25
$ws.Range("K74").Value = 2190.25
$ws.Range("M74").Value = -1316.25

$ws.Range("H77").Value = 2153.2222
$ws.Range("I77").Value = 2190.25
$ws.Range("K77").Value = 10951.25
$ws.Range("M77").Value = -6583.25

$ws.Range("H122").Value = 1911.9524
$ws.Range("I122").Value = 1622.2307
$ws.Range("J122").Value = 2382.75
$ws.Range("K122").Value = 4866.6921
$ws.Range("L122").Value = 7148.25
$ws.Range("M122").Value = -2416.6921
$ws.Range("N122").Value = -12048.25

$ws.Range("H136").Value = 3461.44
$ws.Range("I136").Value = 3461.44
$ws.Range("K136").Value = 10384.32
$ws.Range("M136").Value = -7834.32

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 105.35294
$ws.Range("I4").Value = 42.18182
$ws.Range("J4").Value = 221.16667
$ws.Range("K4").Value = 42.18182
$ws.Range("L4").Value = 221.16667
$ws.Range("M4").Value = 72.81818
$ws.Range("N4").Value = -451.16667

$ws.Range("H94").Value = 3028
$ws.Range("I94").Value = 2447
$ws.Range("K94").Value = 2447
$ws.Range("M94").Value = -1996

$ws.Range("H96").Value = 14974.2
$ws.Range("I96").Value = 7194.1113
$ws.Range("J96").Value = 84995
$ws.Range("K96").Value = 7194.1113
$ws.Range("L96").Value = 84995
$ws.Range("M96").Value = -4448.1113
$ws.Range("N96").Value = -90487

$ws.Range("H122").Value = 70390
$ws.Range("J122").Value = 70390
$ws.Range("L122").Value = 70390
$ws.Range("N122").Value = -80190

$ws.Range("H134").Value = 2615.0527
$ws.Range("I134").Value = 2569.7646
$ws.Range("K134").Value = 7709.293799999999
$ws.Range("M134").Value = -5174.293799999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1295.25
$ws.Range("I58").Value = 1337.4286
$ws.Range("K58").Value = 1337.4286
$ws.Range("M58").Value = -1134.4286

$ws.Range("H136").Value = 1295.25
$ws.Range("I136").Value = 1337.4286
$ws.Range("K136").Value = 4012.2858
$ws.Range("M136").Value = -1462.2858

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 99997
$ws.Range("J37").Value = 99997
$ws.Range("L37").Value = 299991
$ws.Range("N37").Value = -300215

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 12334.333
$ws.Range("I102").Value = 8501.75
$ws.Range("K102").Value = 8501.75
$ws.Range("M102").Value = -6879.75

$ws.Range("H122").Value = 33422.46
$ws.Range("I122").Value = 36449.3
$ws.Range("K122").Value = 109347.9
$ws.Range("M122").Value = -106897.9

$ws.Range("H132").Value = 2974.1304
$ws.Range("I132").Value = 2930.35
$ws.Range("J132").Value = 3266
$ws.Range("K132").Value = 8791.049999999999
$ws.Range("L132").Value = 9798
$ws.Range("M132").Value = -6261.049999999999
$ws.Range("N132").Value = -14858

$ws.Range("H134").Value = 45333.5
$ws.Range("J134").Value = 45333.5
$ws.Range("L134").Value = 136000.5
$ws.Range("N134").Value = -141070.5

$ws.Range("H140").Value = 229999
$ws.Range("J140").Value = 229999
$ws.Range("L140").Value = 229999
$ws.Range("N140").Value = -240359

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6638.7
$ws.Range("I7").Value = 5833
$ws.Range("J7").Value = 7847.25
$ws.Range("K7").Value = 5833
$ws.Range("L7").Value = 7847.25
$ws.Range("M7").Value = -5721
$ws.Range("N7").Value = -8071.25

$ws.Range("H40").Value = 5456.5713
$ws.Range("I40").Value = 3752.923
$ws.Range("J40").Value = 8225
$ws.Range("K40").Value = 3752.923
$ws.Range("L40").Value = 8225
$ws.Range("M40").Value = -3616.923
$ws.Range("N40").Value = -8497

$ws.Range("H126").Value = 6638.7
$ws.Range("I126").Value = 5833
$ws.Range("J126").Value = 7847.25
$ws.Range("K126").Value = 17499
$ws.Range("L126").Value = 23541.75
$ws.Range("M126").Value = -15029
$ws.Range("N126").Value = -28481.75

$ws.Range("H132").Value = 3778.8386
$ws.Range("I132").Value = 3630.261
$ws.Range("J132").Value = 4206
$ws.Range("K132").Value = 10890.783
$ws.Range("L132").Value = 12618
$ws.Range("M132").Value = -8360.782999999999
$ws.Range("N132").Value = -17678

$ws.Range("H136").Value = 14819.762
$ws.Range("I136").Value = 3246.5
$ws.Range("K136").Value = 9739.5
$ws.Range("M136").Value = -7189.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 184999.5
$ws.Range("J64").Value = 184999.5
$ws.Range("L64").Value = 184999.5
$ws.Range("N64").Value = -185495.5

$ws.Range("H67").Value = 184999.5
$ws.Range("J67").Value = 184999.5
$ws.Range("L67").Value = 184999.5
$ws.Range("N67").Value = -186715.5

$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()

$ws.Range("H107").Value = 3171.8333
$ws.Range("I107").Value = 3074.75
$ws.Range("K107").Value = 9224.25
$ws.Range("M107").Value = -7304.25

$ws.Range("H126").Value = 1909.8572
$ws.Range("I126").Value = 1365.3334
$ws.Range("K126").Value = 4096.0002
$ws.Range("M126").Value = -1626.0002

$ws.Range("H132").Value = 1395
$ws.Range("I132").Value = 1353.1936
$ws.Range("K132").Value = 4059.5808
$ws.Range("M132").Value = -1529.5808

$ws.Range("H136").Value = 994.7619
$ws.Range("I136").Value = 1041.6666
$ws.Range("J136").Value = 713.3333
$ws.Range("K136").Value = 3124.9998
$ws.Range("L136").Value = 2139.9999
$ws.Range("M136").Value = -574.9998000000001
$ws.Range("N136").Value = -7239.9999

$ws.Range("H138").Value = 106999.664
$ws.Range("J138").Value = 106999.664
$ws.Range("L138").Value = 106999.664
$ws.Range("N138").Value = -117279.664

$ws.Range("H139").Value = 129234.5
$ws.Range("J139").Value = 128979.336
$ws.Range("L139").Value = 128979.336
$ws.Range("N139").Value = -139259.336
